$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly indexing bug-fix: the per-row qoq-error series was missing its
# leading (quarter-0) observation, so every later observation had slid one
# column to the left of where it belonged. For each data row (2-16) we now
# insert the correct quarter-0 value into column B and push the existing
# values one column to the right (B->C, C->D, ... J->K); anything that would
# overflow past column K (the sheet's last populated column) is discarded,
# matching the original row lengths.

$lastCol = 11  # column K

$newFirstValues = @{
    2  = -1.025188112727922
    3  = 0.08364543516793629
    4  = -0.1538585523806955
    5  = 0.7495351060200912
    6  = 0.03849281619118239
    7  = -0.2590580299438133
    8  = 0.01855976243503714
    9  = 0.1467044301255134
    10 = -0.1819613811903656
    11 = 0.4718454808444464
    12 = -0.08594117411414147
    13 = -0.07695400962807622
    14 = -0.5068991247689255
    15 = 0.6215838649243215
    16 = -0.2766911554241067
}

for ($row = 2; $row -le 16; $row++) {
    for ($c = $lastCol; $c -ge 3; $c--) {
        $srcVal = $ws.Cells.Item($row, $c - 1).Value2
        $ws.Cells.Item($row, $c).Value2 = $srcVal
    }
    $ws.Cells.Item($row, 2).Value2 = $newFirstValues[$row]
}
